# Republishing after code cleanup
# Insert two new config rows into the "Constants" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- Insert new row with "ACMESystem1_WorkItemURL" before the current row 14 ---
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "ACMESystem1_WorkItemURL"
$ws.Range("C14").Value = "URL to open a specific work item in the ACME System 1 portal"
$ws.Range("B14").Value = "https://acme-test.uipath.com/work-items/{0}"

# --- Insert new row with "Process_DefaultDownloadLocation" before the (now shifted) row 17 ---
$ws.Rows.Item(17).Insert()
$ws.Range("C17").Value = "Default path to store the downloaded files"
$ws.Range("B17").Value = "D:\UiPath Project Space\ACMESystem1_ProcessVendorInvoice\Invoices"
$ws.Range("A17").Value = "Process_DefaultDownloadLocation"

# --- Update the view so the selection / scroll position matches the saved file ---
$ws.Activate()
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 6
